# Convention change to support multi-axle vehicles:
# rename the front/rear axle position labels ("sAxleF"/"sAxleR") to the
# generic multi-axle labels ("sAxle1"/"sAxle2") on every vehicle sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sedan_HambaLG", "Sedan_Hamba", "Bus_Makhulu")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("A6").Value = "sAxle2"
}

# Sedan_Hamba: last touched cells were A5:A6, then A17
$ws2 = $wb.Worksheets.Item("Sedan_Hamba")
$ws2.Range("A5:A6").Select()
$ws2.Range("A17").Select()

# Bus_Makhulu: last touched cells were A5:A6, then A2
$ws3 = $wb.Worksheets.Item("Bus_Makhulu")
$ws3.Range("A5:A6").Select()
$ws3.Range("A2").Select()

# Sedan_HambaLG ends up the active/selected tab, with A2 as the final
# active cell.
$ws1 = $wb.Worksheets.Item("Sedan_HambaLG")
$ws1.Activate()
$ws1.Range("A2").Select()
